# Rebuild Sheet1 with the new "Taxi / Doctor appointment" schedule data,
# applying the same look (left-aligned text columns, wrapped Address/Sub
# Category column, wrapped Result column) used on the existing sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Start clean.
$ws.Cells.Clear()

# Header row.
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Number"
$ws.Range("C1").Value = "Address"
$ws.Range("D1").Value = "Sub Category"
$ws.Range("E1").Value = "Result"

# Data rows.
$names    = @("Pawan10", "Ram10", "Pawan11", "Ram11", "Pawan12", "Ram12", "Pawan13", "Ram13")
$numbers  = @(90256985424, 90256985421, 90256985418, 90256985415, 90256985412, 90256985409, 90256985406, 90256985403)
$address  = @("Madhapur", "Hyderabad", "Madhapur", "Hyderabad", "Madhapur", "Hyderabad", "Madhapur", "Hyderabad")
$subcat   = @("Taxi", "Doctor appointment", "Taxi", "Doctor appointment", "Taxi", "Doctor appointment", "Taxi", "Doctor appointment")
$results  = @("Pass", "Pass", "Pass", "Pass", "Pass", "Pass", "Pass", $null)

for ($i = 0; $i -lt $names.Length; $i++) {
    $r = $i + 2

    $ws.Cells.Item($r, 1).Value = $names[$i]
    $ws.Cells.Item($r, 2).Value = $numbers[$i]
    $ws.Cells.Item($r, 3).Value = $address[$i]
    $ws.Cells.Item($r, 4).Value = $subcat[$i]

    if ($results[$i] -ne $null) {
        $ws.Cells.Item($r, 5).Value = $results[$i]
    }
}

# Formatting: columns A:C left-aligned, column D left-aligned + wrapped,
# column E (Result) wrapped.
$ws.Range("A1:C9").HorizontalAlignment = -4131
$ws.Range("D1:D9").HorizontalAlignment = -4131
$ws.Range("D1:D9").WrapText = $true
$ws.Range("E2:E8").WrapText = $true

# Column widths, matching the other sheets' look (best-fit-style widths).
$ws.Columns.Item(1).ColumnWidth = 8
$ws.Columns.Item(2).ColumnWidth = 11.14
$ws.Columns.Item(3).ColumnWidth = 9.6
$ws.Columns.Item(4).ColumnWidth = 22
$ws.Columns.Item(5).ColumnWidth = 5.6

# Touch Sheet2's trailing formatted-but-empty cell so its style re-resolves
# through the (now reorganised) style table, same as Sheet1's.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("D4").WrapText = $true

$ws.Range("E11").Select()
